$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 2.18
$ws.Range("I2").Value = 2.24
$ws.Range("O2").Value = 1.37
$ws.Range("U2").Value = 2.08

# Row 3
$ws.Range("N3").Value = 1.26
$ws.Range("O3").Value = 1.2
$ws.Range("P3").Value = 1.25
$ws.Range("Q3").Value = 1.2
$ws.Range("S3").Value = 1.2

# Row 7
$ws.Range("R7").Value = 1.57
$ws.Range("S7").Value = 2.02

# Row 8
$ws.Range("G8").Value = 2.5
$ws.Range("I8").Value = 3.85

# Row 9
$ws.Range("F9").Value = 4.7
$ws.Range("G9").Value = 7.6
$ws.Range("H9").Value = 1.6
$ws.Range("I9").Value = 1.75
$ws.Range("J9").Value = 4
$ws.Range("K9").Value = 4.6
$ws.Range("P9").Value = 2.02
$ws.Range("Q9").Value = 1.72

$wb.Save()
